# Loan RBI, Variable Instalments
#
# The "Repayment Schedule" sheet gains a new (blank) column inserted just
# before the existing "Late" column, pushing "Late" / "Heading" / "Outstanding"
# one column to the right (N->O, O->P, P->Q). We reproduce that with a
# standard Excel "insert entire column" operation, which shifts all the
# cells/values/styles at and after column N one column to the right and
# leaves the freshly inserted column N blank.
#
# The workbook also ends up with the "Repayment Schedule" tab active/selected
# (with the cursor on L15) instead of the "Transactions" tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N (shifts N:P -> O:Q).
$ws.Range("N1").EntireColumn.Insert() | Out-Null

# Make "Repayment Schedule" the active sheet/tab with L15 selected.
$ws.Activate() | Out-Null
$ws.Range("L15").Select() | Out-Null
